$wb = $excel.ActiveWorkbook

# --- Add the new "Estimates" worksheet, placed after the existing Sheet1 ---
$sheetCount = $wb.Worksheets.Count
$ws2 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($sheetCount))
$ws2.Name = "Estimates"

# --- Column widths (B ~56 chars, C ~10.5 chars) ---
$ws2.Columns.Item(2).ColumnWidth = 55.15
$ws2.Columns.Item(3).ColumnWidth = 9.65

# --- Content ---
$ws2.Range("B2").Value = "Total Story Points Esimates (Including Desing, Cut Effort, DB Design, Testing, Requirement Detailing, Code Review, Bug Fixing, Documentation, Release Notes)"
$ws2.Range("B2").WrapText = $true
$ws2.Range("C2").Value = 314

# Row 2 needs to be tall enough to show the wrapped text
$ws2.Rows.Item(2).RowHeight = 43.5

# Page setup (portrait, matching the rest of the workbook)
$ws2.PageSetup.Orientation = 1

# Leave the selection on C2, matching the authored file, and make
# "Estimates" the active (visible) sheet/tab.
$ws2.Range("C2").Select() | Out-Null
$ws2.Activate() | Out-Null
